$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, reusing the same formatting (bold, centered,
# bordered) as the other header cells (e.g. G1) by copying its format.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data values for the new "Save" column (rows 2-8), plain numbers with
# the default (unstyled) cell format like the rest of the data rows.
$saveValues = @(0, 1, 0, 0, 0, 1, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
